$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Extend Table1 to cover the 5 new rows (A1:F8 -> A1:F10)
$lo.Resize($ws.Range("A1:F10"))

# --- Row 6: OLED Display ---
$ws.Range("B6").Value2 = "OLED Display"
$urlF6 = "https://www.kiwi-electronics.com/nl/monochroom-0-96quot-128x64-oled-grafisch-display-754?search=oled"
$ws.Hyperlinks.Add($ws.Range("F6"), $urlF6, [Type]::Missing, [Type]::Missing, $urlF6) | Out-Null
$ws.Range("F6").Value2 = "Monochroom 0.96"" 128x64 OLED grafisch display (kiwi-electronics.com)"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("C6").Value2 = 23.5
$ws.Range("D6").Value2 = 1
$ws.Range("E6").Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

# --- Descriptions for rows 7-9 ---
$ws.Range("B7").Value2 = "N-MOSFET"
$ws.Range("B8").Value2 = "BJT"
$ws.Range("B9").Value2 = "USB-C connector"

# --- Row 9: USB-C connector ---
$urlF9 = "https://www.lcsc.com/product-detail/USB-Connectors_XKB-Connectivity-U262-161N-4BVC11_C319148.html"
$ws.Hyperlinks.Add($ws.Range("F9"), $urlF9, [Type]::Missing, [Type]::Missing, $urlF9) | Out-Null
$ws.Range("F9").Value2 = "U262-161N-4BVC11 XKB Connectivity | C319148 - LCSC Electronics"
$ws.Range("F9").Style = "Hyperlink"

# --- Row 7: N-MOSFET ---
$urlF7 = "https://www.lcsc.com/product-detail/MOSFETs_Diodes-Incorporated-DMN2056U-7_C332302.html"
$ws.Hyperlinks.Add($ws.Range("F7"), $urlF7, [Type]::Missing, [Type]::Missing, $urlF7) | Out-Null
$ws.Range("F7").Value2 = "DMN2056U-7 Diodes Incorporated | C332302 - LCSC Electronics"
$ws.Range("F7").Style = "Hyperlink"
$ws.Range("A7").Value2 = "DMN2056U-7"

# --- Row 8: BJT ---
$urlF8 = "https://www.lcsc.com/product-detail/Bipolar-Transistors-BJT_onsemi-SBC847BDW1T1G_C232475.html"
$ws.Hyperlinks.Add($ws.Range("F8"), $urlF8, [Type]::Missing, [Type]::Missing, $urlF8) | Out-Null
$ws.Range("F8").Value2 = "SBC847BDW1T1G onsemi | C232475 - LCSC Electronics"
$ws.Range("F8").Style = "Hyperlink"
$ws.Range("A8").Value2 = "SBC847BDW1T1G"

# --- Row 9 Part ---
$ws.Range("A9").Value2 = "U262-161N-4BVC11"

# --- Row 10: LDO ---
$ws.Range("B10").Value2 = "LDO"
$urlF10 = "https://www.lcsc.com/product-detail/Linear-Voltage-Regulators-LDO_Diodes-Incorporated-AZ1117H-3-3TRE1_C92517.html"
$ws.Hyperlinks.Add($ws.Range("F10"), $urlF10, [Type]::Missing, [Type]::Missing, $urlF10) | Out-Null
$ws.Range("F10").Value2 = "AZ1117H-3.3TRE1 Diodes Incorporated | C92517 - LCSC Electronics"
$ws.Range("F10").Style = "Hyperlink"
$ws.Range("A10").Value2 = "AZ1117H-3.3TRE1"

# --- Remaining numeric cells + formulas for rows 7-10 ---
$ws.Range("C7").Value2 = 0.62
$ws.Range("D7").Value2 = 2
$ws.Range("E7").Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

$ws.Range("C8").Value2 = 0.62
$ws.Range("D8").Value2 = 1
$ws.Range("E8").Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

$ws.Range("C9").Value2 = 0.39
$ws.Range("D9").Value2 = 1
$ws.Range("E9").Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

$ws.Range("C10").Value2 = 0.12
$ws.Range("D10").Value2 = 1
$ws.Range("E10").Formula = "=Table1[[#This Row],[Price]]*Table1[[#This Row],[Qty]]"

# Column F got wider to fit the new, longer descriptions
$ws.Columns("F").ColumnWidth = 73.9

# Selection moved down past the new rows
$ws.Range("A12").Select() | Out-Null
